# Update cryptocurrency price and volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.624.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.460.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +3.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.459.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.164"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.527.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "337.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0821"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.87%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "425.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "130.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.483"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.563"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
